$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '42.211.42'
$ws.Cells.Item(2, 5).Value = '  -1.87%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.249.67'
$ws.Cells.Item(3, 5).Value = '  -2.11%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '247.14'
$ws.Cells.Item(5, 5).Value = '  -2.07%  '
$ws.Cells.Item(6, 5).Value = '  -1.31%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '77.18'
$ws.Cells.Item(7, 5).Value = '  +3.87%  '
$ws.Cells.Item(8, 5).Value = '  -0.01%  '
$ws.Cells.Item(9, 5).Value = '  -4.53%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '42.13'
$ws.Cells.Item(10, 5).Value = '  +5.93%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0953'
$ws.Cells.Item(11, 5).Value = '  -2.93%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '7.12'
$ws.Cells.Item(12, 5).Value = '  -4.88%  '
$ws.Cells.Item(13, 5).Value = '  -3.09%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '2.585.07'
$ws.Cells.Item(14, 5).Value = '  -2.17%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '14.76'
$ws.Cells.Item(15, 5).Value = '  -4.25%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '0.861'
$ws.Cells.Item(16, 5).Value = '  -1.77%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '2.239.21'
$ws.Cells.Item(17, 5).Value = '  -3.64%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '42.092.64'
$ws.Cells.Item(18, 5).Value = '  -1.98%  '
$ws.Cells.Item(19, 5).Value = '  -2.86%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '72.00'
$ws.Cells.Item(20, 5).Value = '  -1.09%  '
$ws.Cells.Item(21, 5).Value = '  -3.09%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '2.29'
$ws.Cells.Item(22, 5).Value = '  +0.70%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '232.04'
$ws.Cells.Item(23, 5).Value = '  -2.75%  '
$ws.Cells.Item(24, 5).Value = '  -0.06%  '
$ws.Cells.Item(25, 5).Value = '  -2.75%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '3.61'
$ws.Cells.Item(26, 5).Value = '  -7.52%  '
$ws.Cells.Item(27, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '7.57'
$ws.Cells.Item(27, 5).Value = '  +18.79%  '
$ws.Cells.Item(28, 2).Value = 'PancakeSwap'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.29'
$ws.Cells.Item(28, 5).Value = '  -5.41%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '170.14'
$ws.Cells.Item(29, 5).Value = '  +1.58%  '
$ws.Cells.Item(30, 5).Value = '  -2.12%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '20.63'
$ws.Cells.Item(31, 5).Value = '  -2.42%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '0.0838'
$ws.Cells.Item(32, 5).Value = '  -0.39%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '32.71'
$ws.Cells.Item(33, 5).Value = '  +5.47%  '
$ws.Cells.Item(34, 5).Value = '  -5.31%  '
$ws.Cells.Item(35, 5).Value = '  -1.40%  '
$ws.Cells.Item(36, 5).Value = '  -2.51%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '4.94'
$ws.Cells.Item(37, 5).Value = '  +2.47%  '
$ws.Cells.Item(38, 2).Value = 'Celestia'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '14.39'
$ws.Cells.Item(38, 5).Value = '  +3.62%  '
$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.0303'
$ws.Cells.Item(39, 5).Value = '  -2.24%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '5.90'
$ws.Cells.Item(40, 5).Value = '  -0.14%  '
$ws.Cells.Item(41, 5).Value = '  -7.49%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '112.73'
$ws.Cells.Item(42, 5).Value = '  +7.32%  '
$ws.Cells.Item(43, 5).Value = '  -7.21%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '61.03'
$ws.Cells.Item(44, 5).Value = '  -2.13%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '8.69'
$ws.Cells.Item(45, 5).Value = '  -5.47%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.0993'
$ws.Cells.Item(46, 5).Value = '  -4.17%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.998'
$ws.Cells.Item(47, 5).Value = '  -0.58%  '
$ws.Cells.Item(48, 5).Value = '  -4.23%  '
$ws.Cells.Item(49, 5).Value = '  -1.64%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '4.30'
$ws.Cells.Item(50, 5).Value = '  -12.29%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.441'
$ws.Cells.Item(51, 5).Value = '  +15.15%  '
